# Update the Global Glider Cal and Ingest sheet (Asset_Cal_Info):
#  - CC_scattering_angle value (row 2, col F): 117 -> 140
#  - CC_angular_resolution value (row 4, col F): 1.08 -> 1.13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

$ws.Range("F2").Value = 140
$ws.Range("F4").Value = 1.13

# Leave the selection on the last-edited cell, matching the saved view state.
$ws.Range("F4").Select() | Out-Null
